# Update the "Generate" timestamps on the handback-status report.
# Sheet1 = "Overview", Sheet2 = "zh-cn", Sheet3 = "de-de"
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# "Latest HO Xliff Generate Date" for the de-de handoff file.
# Overview!G2 and de-de!H2 share the same original value, so both move
# together to stay in sync.
$wsOverview.Range("G2").Value = "2016-09-06 11:08:56"
$wsDeDe.Range("H2").Value     = "2016-09-06 11:08:56"

# zh-cn sheet: Correspond Handoff / Handback Datetime for the first row.
$wsZhCn.Range("H2").Value = "2016-09-06 11:08:45"
$wsZhCn.Range("K2").Value = "2016-09-06 11:09:36"

# de-de sheet: Correspond Handback Datetime for the first row.
$wsDeDe.Range("K2").Value = "2016-09-06 11:09:55"
